$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ranames = @(
    'Michael Walton',
    'Samuel Lester',
    'Jeffrey Campos',
    'Christopher Garrett',
    'Melissa Welch',
    'Dr. Andre Craig',
    'Jacob Hutchinson',
    'Alicia Graham',
    'Kaitlin Williamson',
    'Beth Villanueva',
    'Andrew Roman',
    'Scott Orozco',
    'Joseph Lee',
    'Eileen Hill',
    'Victor Fitzgerald',
    'Raymond Fuller',
    'April Flores',
    'Robert Wright',
    'James Becker',
    'Travis Villarreal',
    'Dylan Schultz',
    'Adam Bennett',
    'Erica Thomas',
    'Eric Cobb',
    'Michael Briggs',
    'Richard Clark',
    'Michael Taylor',
    'Sabrina Hancock',
    'Jill Travis',
    'Breanna Jensen',
    'Karen Bowman',
    'Karen Webb',
    'Joshua Cox',
    'Thomas Hester',
    'Sandra Hobbs',
    'Kevin Lee',
    'Deanna Blair',
    'Michelle Williams',
    'Anthony Jones',
    'David Jensen',
    'Kenneth Mullen',
    'Matthew Robinson',
    'Ashley Hopkins',
    'Linda Smith',
    'Aaron Nichols',
    'Lauren Tyler',
    'Jeffrey Carlson',
    'Andre Howard',
    'Stephen Suarez',
    'Robin Santiago',
    'Andrew Sullivan',
    'Maria Meza',
    'Roberta Jenkins',
    'Meghan Dunn',
    'Desiree Brock',
    'Laura Watson',
    'Valerie Cohen',
    'Jacqueline Buchanan',
    'Kimberly Christensen',
    'Brian Jones',
    'Melinda Ramirez',
    'Monica Olsen',
    'Caitlin Flores',
    'Wendy Waters',
    'Kristina Torres',
    'Crystal Nielsen',
    'Eric Chen',
    'Sarah Jordan',
    'Austin Solomon',
    'Rose Mitchell',
    'Juan Rodriguez',
    'Jennifer Gibson',
    'Linda Miller',
    'Kenneth Owens',
    'Victor Martinez',
    'Andrew Baldwin',
    'Andrea Bryan',
    'Kerry Day',
    'Daniel Carpenter',
    'Jesse Flores',
    'Michael Grant',
    'Michael Simmons',
    'Scott Morales',
    'Jose Medina',
    'Nina Murphy',
    'Kenneth Higgins',
    'Brittany Nelson',
    'Dr. Patricia Gill',
    'Sabrina Patterson',
    'Morgan Hernandez',
    'Whitney Powell',
    'Thomas Spencer',
    'Garrett Williams',
    'Jill Jackson',
    'Bradley Sullivan',
    'Jessica Moore',
    'Mary Hernandez',
    'Christopher Castaneda',
    'Danielle Wilkins',
    'Cindy Pierce',
    'Steve Cooper',
    'William Brown',
    'Jessica Bradshaw',
    'Christina Jones',
    'Maurice Black',
    'Sarah Garcia',
    'Christian Richardson',
    'Dr. Norma Ramirez MD',
    'Dominique Valdez',
    'Mallory Logan',
    'Robert Mills',
    'Gina Rios',
    'Laura Frazier',
    'Crystal Lopez',
    'Deanna Cabrera',
    'Amanda Campbell',
    'Mary Reyes',
    'Stephanie Collins',
    'Matthew Doyle',
    'Julie Ewing',
    'Deborah Lang',
    'Brian Spears',
    'Ryan Bennett',
    'Valerie Sanders',
    'Justin Spencer',
    'Lisa Mcbride',
    'Mary Mckinney',
    'Victoria Malone',
    'James Castillo',
    'Aaron Hawkins',
    'Charles Harris',
    'Kyle Conway',
    'Michelle Davis',
    'Sean Russell',
    'Leslie Callahan',
    'Monica Wallace',
    'Bruce English',
    'April Dawson',
    'Kimberly Berger',
    'George Banks',
    'Sheila Mendoza',
    'Isabella Johnson',
    'Kathleen Gonzales',
    'Gabriela Jackson',
    'Victoria Frederick',
    'Krystal Kerr',
    'Angela Velez',
    'Jacqueline Gonzalez',
    'Edward Conway',
    'Kristine Smith',
    'Katherine Cole',
    'David Wang',
    'Paul Walter',
    'Brenda Thompson',
    'Ethan Tucker',
    'David Jackson',
    'Denise Carlson',
    'Christopher Flores Jr.',
    'Paula Hanson',
    'Lisa Russell'
)

$teamids = @(
    '9cba5ad7-314e-4f2a-80a9-fc31cbf3f0c7',
    '61f0969e-22a4-4374-8588-d6511915b05e',
    '5e1a20f6-82bf-4dee-aa79-41702d9feb41',
    'e5c40d19-b03a-4f5a-82c8-25540cd45e07',
    'e5c40d19-b03a-4f5a-82c8-25540cd45e07',
    '5e1a20f6-82bf-4dee-aa79-41702d9feb41',
    '61f0969e-22a4-4374-8588-d6511915b05e',
    '630f61e8-543f-46e2-af63-2b62e8bc4fd2',
    '94742748-e7ab-454b-8ff6-9893440bd059',
    '9cba5ad7-314e-4f2a-80a9-fc31cbf3f0c7',
    '94742748-e7ab-454b-8ff6-9893440bd059',
    'e5ca6e2b-5f54-4acd-ad7b-03e631313986',
    '9cba5ad7-314e-4f2a-80a9-fc31cbf3f0c7',
    'e0228b4f-7807-45db-a3f6-8c6e1f4adf41',
    '6afc31f0-3916-443a-92c4-b5eb425a9bc3',
    'e0228b4f-7807-45db-a3f6-8c6e1f4adf41',
    'e5c40d19-b03a-4f5a-82c8-25540cd45e07',
    'e0228b4f-7807-45db-a3f6-8c6e1f4adf41',
    '0eeb011c-24fb-4476-91f7-d8e28ae49c2f',
    '9cba5ad7-314e-4f2a-80a9-fc31cbf3f0c7',
    '630f61e8-543f-46e2-af63-2b62e8bc4fd2',
    '6afc31f0-3916-443a-92c4-b5eb425a9bc3',
    'e5c40d19-b03a-4f5a-82c8-25540cd45e07',
    '0435a227-38e8-494e-b1bf-271b00893eae',
    '94742748-e7ab-454b-8ff6-9893440bd059',
    'e5ca6e2b-5f54-4acd-ad7b-03e631313986',
    'e0228b4f-7807-45db-a3f6-8c6e1f4adf41',
    'e0228b4f-7807-45db-a3f6-8c6e1f4adf41',
    '94742748-e7ab-454b-8ff6-9893440bd059',
    '5e1a20f6-82bf-4dee-aa79-41702d9feb41',
    '5e1a20f6-82bf-4dee-aa79-41702d9feb41',
    '7fc75193-58a0-4e7d-ab42-382ec10a8be4',
    '6afc31f0-3916-443a-92c4-b5eb425a9bc3',
    '9cba5ad7-314e-4f2a-80a9-fc31cbf3f0c7',
    'e5c40d19-b03a-4f5a-82c8-25540cd45e07',
    '94742748-e7ab-454b-8ff6-9893440bd059',
    '94742748-e7ab-454b-8ff6-9893440bd059',
    '9cba5ad7-314e-4f2a-80a9-fc31cbf3f0c7',
    '630f61e8-543f-46e2-af63-2b62e8bc4fd2',
    '6afc31f0-3916-443a-92c4-b5eb425a9bc3',
    'e5c40d19-b03a-4f5a-82c8-25540cd45e07',
    'e5ca6e2b-5f54-4acd-ad7b-03e631313986',
    'e5ca6e2b-5f54-4acd-ad7b-03e631313986',
    '630f61e8-543f-46e2-af63-2b62e8bc4fd2',
    '7fc75193-58a0-4e7d-ab42-382ec10a8be4',
    '61f0969e-22a4-4374-8588-d6511915b05e',
    '7fc75193-58a0-4e7d-ab42-382ec10a8be4',
    '7fc75193-58a0-4e7d-ab42-382ec10a8be4',
    '6afc31f0-3916-443a-92c4-b5eb425a9bc3',
    '9cba5ad7-314e-4f2a-80a9-fc31cbf3f0c7',
    '0eeb011c-24fb-4476-91f7-d8e28ae49c2f',
    '6afc31f0-3916-443a-92c4-b5eb425a9bc3',
    '61f0969e-22a4-4374-8588-d6511915b05e',
    '9cba5ad7-314e-4f2a-80a9-fc31cbf3f0c7',
    '9cba5ad7-314e-4f2a-80a9-fc31cbf3f0c7',
    '61f0969e-22a4-4374-8588-d6511915b05e',
    'e0228b4f-7807-45db-a3f6-8c6e1f4adf41',
    '7fc75193-58a0-4e7d-ab42-382ec10a8be4',
    'e5ca6e2b-5f54-4acd-ad7b-03e631313986',
    '61f0969e-22a4-4374-8588-d6511915b05e',
    '7fc75193-58a0-4e7d-ab42-382ec10a8be4',
    '61f0969e-22a4-4374-8588-d6511915b05e',
    'e5ca6e2b-5f54-4acd-ad7b-03e631313986',
    '7fc75193-58a0-4e7d-ab42-382ec10a8be4',
    '61f0969e-22a4-4374-8588-d6511915b05e',
    '9cba5ad7-314e-4f2a-80a9-fc31cbf3f0c7',
    '0435a227-38e8-494e-b1bf-271b00893eae',
    '5e1a20f6-82bf-4dee-aa79-41702d9feb41',
    '630f61e8-543f-46e2-af63-2b62e8bc4fd2',
    '0435a227-38e8-494e-b1bf-271b00893eae',
    'e5c40d19-b03a-4f5a-82c8-25540cd45e07',
    '7fc75193-58a0-4e7d-ab42-382ec10a8be4',
    '5e1a20f6-82bf-4dee-aa79-41702d9feb41',
    '6afc31f0-3916-443a-92c4-b5eb425a9bc3',
    '61f0969e-22a4-4374-8588-d6511915b05e',
    '5e1a20f6-82bf-4dee-aa79-41702d9feb41',
    '0eeb011c-24fb-4476-91f7-d8e28ae49c2f',
    '94742748-e7ab-454b-8ff6-9893440bd059',
    '94742748-e7ab-454b-8ff6-9893440bd059',
    '7fc75193-58a0-4e7d-ab42-382ec10a8be4',
    '5184566d-523a-4432-848d-ac234ffb6ac6',
    'de503c24-f17d-47a9-9a47-6f0a194f8c9c',
    '4279fd55-c2c1-440d-abaa-430f3c27be44',
    '3b8adc57-0f6f-482c-8306-9830e819d666',
    '1ad987f9-6aa5-4e4d-8f1a-e8bea8fa4fcd',
    'e7bb31c1-e095-453b-95ff-565ea62efb0a',
    '7ba6e30b-04d4-4449-b7a9-2a6c7bb23764',
    '4279fd55-c2c1-440d-abaa-430f3c27be44',
    '4279fd55-c2c1-440d-abaa-430f3c27be44',
    '718c6b8f-7c00-4bcb-b53c-8f3f42154362',
    '1ad987f9-6aa5-4e4d-8f1a-e8bea8fa4fcd',
    '7ba6e30b-04d4-4449-b7a9-2a6c7bb23764',
    '3b8adc57-0f6f-482c-8306-9830e819d666',
    '850a92da-c3d6-4fb9-a510-99626e9ad312',
    'e42288a3-b5af-4464-bc45-85d438bcea11',
    '718c6b8f-7c00-4bcb-b53c-8f3f42154362',
    'e7bb31c1-e095-453b-95ff-565ea62efb0a',
    'e7bb31c1-e095-453b-95ff-565ea62efb0a',
    '4279fd55-c2c1-440d-abaa-430f3c27be44',
    '5184566d-523a-4432-848d-ac234ffb6ac6',
    '718c6b8f-7c00-4bcb-b53c-8f3f42154362',
    'de503c24-f17d-47a9-9a47-6f0a194f8c9c',
    '1ad987f9-6aa5-4e4d-8f1a-e8bea8fa4fcd',
    '5184566d-523a-4432-848d-ac234ffb6ac6',
    'de503c24-f17d-47a9-9a47-6f0a194f8c9c',
    '1ad987f9-6aa5-4e4d-8f1a-e8bea8fa4fcd',
    '3b8adc57-0f6f-482c-8306-9830e819d666',
    'e7bb31c1-e095-453b-95ff-565ea62efb0a',
    'de503c24-f17d-47a9-9a47-6f0a194f8c9c',
    '5184566d-523a-4432-848d-ac234ffb6ac6',
    '7ba6e30b-04d4-4449-b7a9-2a6c7bb23764',
    '3b8adc57-0f6f-482c-8306-9830e819d666',
    '57ebee16-96d2-46a6-ab16-2476b305fd91',
    '1ad987f9-6aa5-4e4d-8f1a-e8bea8fa4fcd',
    '1ad987f9-6aa5-4e4d-8f1a-e8bea8fa4fcd',
    'e42288a3-b5af-4464-bc45-85d438bcea11',
    '5184566d-523a-4432-848d-ac234ffb6ac6',
    '5184566d-523a-4432-848d-ac234ffb6ac6',
    '7ba6e30b-04d4-4449-b7a9-2a6c7bb23764',
    '1ad987f9-6aa5-4e4d-8f1a-e8bea8fa4fcd',
    '850a92da-c3d6-4fb9-a510-99626e9ad312',
    'e7bb31c1-e095-453b-95ff-565ea62efb0a',
    'e7bb31c1-e095-453b-95ff-565ea62efb0a',
    'e42288a3-b5af-4464-bc45-85d438bcea11',
    'e42288a3-b5af-4464-bc45-85d438bcea11',
    '4279fd55-c2c1-440d-abaa-430f3c27be44',
    '850a92da-c3d6-4fb9-a510-99626e9ad312',
    '3b8adc57-0f6f-482c-8306-9830e819d666',
    'de503c24-f17d-47a9-9a47-6f0a194f8c9c',
    '718c6b8f-7c00-4bcb-b53c-8f3f42154362',
    '1ad987f9-6aa5-4e4d-8f1a-e8bea8fa4fcd',
    '4279fd55-c2c1-440d-abaa-430f3c27be44',
    '3b8adc57-0f6f-482c-8306-9830e819d666',
    '1ad987f9-6aa5-4e4d-8f1a-e8bea8fa4fcd',
    'e7bb31c1-e095-453b-95ff-565ea62efb0a',
    '850a92da-c3d6-4fb9-a510-99626e9ad312',
    '1ad987f9-6aa5-4e4d-8f1a-e8bea8fa4fcd',
    '3b8adc57-0f6f-482c-8306-9830e819d666',
    '4279fd55-c2c1-440d-abaa-430f3c27be44',
    '7ba6e30b-04d4-4449-b7a9-2a6c7bb23764',
    '3b8adc57-0f6f-482c-8306-9830e819d666',
    '57ebee16-96d2-46a6-ab16-2476b305fd91',
    'e42288a3-b5af-4464-bc45-85d438bcea11',
    '850a92da-c3d6-4fb9-a510-99626e9ad312',
    '4279fd55-c2c1-440d-abaa-430f3c27be44',
    '1ad987f9-6aa5-4e4d-8f1a-e8bea8fa4fcd',
    '3b8adc57-0f6f-482c-8306-9830e819d666',
    '7ba6e30b-04d4-4449-b7a9-2a6c7bb23764',
    '4279fd55-c2c1-440d-abaa-430f3c27be44',
    '1ad987f9-6aa5-4e4d-8f1a-e8bea8fa4fcd',
    'de503c24-f17d-47a9-9a47-6f0a194f8c9c',
    '850a92da-c3d6-4fb9-a510-99626e9ad312',
    '1ad987f9-6aa5-4e4d-8f1a-e8bea8fa4fcd',
    '718c6b8f-7c00-4bcb-b53c-8f3f42154362',
    '3b8adc57-0f6f-482c-8306-9830e819d666',
    '850a92da-c3d6-4fb9-a510-99626e9ad312',
    '57ebee16-96d2-46a6-ab16-2476b305fd91',
    '3b8adc57-0f6f-482c-8306-9830e819d666',
    '1ad987f9-6aa5-4e4d-8f1a-e8bea8fa4fcd',
    'e7bb31c1-e095-453b-95ff-565ea62efb0a'
)

for ($i = 0; $i -lt $ranames.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 13).Value = $ranames[$i]
    $ws.Cells.Item($r, 15).Value = $teamids[$i]
}

